$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Who" for Advances bidding system (row 11) in Sprint 2
$ws.Range("E11").Value = "Oli"

# Add a new Sprint 2 goal: Pictures
$ws.Range("A12").Value = "Pictures"
$ws.Range("D12").Value = "x"
$ws.Range("E12").Value = "Urs"

# Update "Who" for Search (row 9) in Sprint 2
$ws.Range("E9").Value = "Urs/Arun"

# Move selection to the newly added row
$ws.Range("A12").Select()
